# SortSequenceDiagram.pptx edit script
#
# Implements:
#  - "limitIsPresent" -> "isLimitPresent" text fix
#  - "ParseOrder(order)" -> "parseOrder(args)" (split into (,args,))
#  - "ParseLimit(limit)" -> "parseLimit(args)" (split into (,args,))
#  - removal of a stray empty "Rectangle 155" shape
#  - repositioning / resizing of a number of shapes around the
#    ParserUtil / SortCommand / updateFilteredRestaurantList lifelines
#
# Shape.Left/Top/Width/Height are Single-precision (f32) point values in the
# PowerPoint object model and this host truncates (rather than rounds) when
# converting back to EMU on save, so literal point constants below were
# chosen (via exact f32 bit-search) so that point*12700 truncates to the
# exact target EMU value from the source OOXML.

function Get-ShapeById {
    param($shapes, $id)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1. "[limitIsPresent]" -> "[isLimitPresent]"  (TextBox 123, id=124)
# ---------------------------------------------------------------------------
$shp124 = Get-ShapeById $s.Shapes 124
$tr124 = $shp124.TextFrame.TextRange
# full text is "[limitIsPresent]"; "limitIsPresent" occupies chars 2..15
$tr124.Characters(2, 14).Text = "isLimitPresent"

# ---------------------------------------------------------------------------
# 2. "ParseOrder(order)" -> "parseOrder(args)" with (, args, ) as separate
#    runs (TextBox 150, id=151)
# ---------------------------------------------------------------------------
$shp151 = Get-ShapeById $s.Shapes 151
$tr151 = $shp151.TextFrame.TextRange
$tr151.Characters(1, 10).Text = "parseOrder"
$tr151.Characters(12, 5).Text = "args"

# ---------------------------------------------------------------------------
# 3. "ParseLimit(limit)" -> "parseLimit(args)" with (, args, ) as separate
#    runs (TextBox 151, id=152)
# ---------------------------------------------------------------------------
$shp152 = Get-ShapeById $s.Shapes 152
$tr152 = $shp152.TextFrame.TextRange
$tr152.Characters(1, 10).Text = "parseLimit"
$tr152.Characters(12, 5).Text = "args"

# ---------------------------------------------------------------------------
# 4. Remove the stray empty "Rectangle 155" shape (id=156)
# ---------------------------------------------------------------------------
$shp156 = Get-ShapeById $s.Shapes 156
$shp156.Delete()

# ---------------------------------------------------------------------------
# 5. Reposition remaining shapes (EMU targets reached via exact f32 point
#    literals; see header comment).
# ---------------------------------------------------------------------------

# Rectangle 62 (id=128) - "sc:Sort" activation box: move up slightly
$shp128 = Get-ShapeById $s.Shapes 128
$shp128.Top = 208.4832305908203

# Rectangle 62 (id=158) - ":ParserUtil" activation box: move down slightly
$shp158 = Get-ShapeById $s.Shapes 158
$shp158.Top = 69.33063507080078

# TextBox 178 (id=179) - "SortCommand(order, limit)" label
$shp179 = Get-ShapeById $s.Shapes 179
$shp179.Top = 216.39505004882812

# Straight Arrow Connector 179 (id=180)
$shp180 = Get-ShapeById $s.Shapes 180
$shp180.Top = 235.00717163085938

# Straight Arrow Connector 180 (id=181)
$shp181 = Get-ShapeById $s.Shapes 181
$shp181.Top = 258.0

# Rectangle 188 (id=189)
$shp189 = Get-ShapeById $s.Shapes 189
$shp189.Top = 242.15426635742188

# Elbow Connector 224 (id=225) - keeps rot="16200000" flipH="1"
$shp225 = Get-ShapeById $s.Shapes 225
$shp225.Left = 544.3555297851562
$shp225.Top = 441.1547546386719
$shp225.Width = 7.859606742858887
$shp225.Height = 4.989921569824219

# Rectangle 225 (id=226) - only size changes
$shp226 = Get-ShapeById $s.Shapes 226
$shp226.Width = 8.783543586730957
$shp226.Height = 38.71787643432617

# TextBox 236 (id=237) - "updateFilteredRestaurantList()" label
$shp237 = Get-ShapeById $s.Shapes 237
$shp237.Left = 550.8986206054688
$shp237.Top = 412.66986083984375

# Straight Arrow Connector 237 (id=238)
$shp238 = Get-ShapeById $s.Shapes 238
$shp238.Left = 556.131103515625
$shp238.Top = 423.5941162109375
$shp238.Width = 156.11544799804688
$shp238.Height = 0.0

# Straight Arrow Connector 251 (id=252) - also gains a vertical flip
$shp252 = Get-ShapeById $s.Shapes 252
$shp252.VerticalFlip = -1
$shp252.Left = 556.131103515625
$shp252.Top = 431.3874816894531
$shp252.Width = 156.52536010742188
$shp252.Height = 0.16669292747974396
